$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; this shifts rows 18:55 down to 19:56,
# preserving all of their existing data/formatting.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C18").Value = "Metropolitana"
$ws.Range("D18").Value = 45037
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 100112010
$ws.Range("G18").Value = "Achicoria"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 7000
$ws.Range("N18").Value = "`$/caja 16 unidades"
$ws.Range("O18").Value = "Provincia de Quillota"
$ws.Range("P18").Value = 438
$ws.Range("Q18").Value = 16
$ws.Range("R18").Value = "Hortaliza"
